# Better handle copying paragraph styles.
#
# Each Scripture passage block in this document ends with a pair of
# "placeholder" paragraphs (an elided-verses marker "[...]" followed by a
# blank paragraph) that sit between the "MSCJoin" paragraph that precedes
# them and the "MSCParagraph" paragraph that follows them. Those two
# placeholder paragraphs were not getting the "MSCJoin" paragraph style
# copied onto them. Fix that by applying the "MSC_Join" style to the
# "[...]" paragraph and to the blank paragraph immediately after it,
# wherever that pattern occurs.

$d = $word.ActiveDocument
$joinStyleName = "MSC_Join"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($null -eq $text) {
        continue
    }

    $trimmed = $text.Trim([char]13, [char]7)

    if ($trimmed -eq "[...]") {
        $para.Style = $joinStyleName

        if ($i + 1 -le $count) {
            $nextPara = $d.Paragraphs.Item($i + 1)
            $nextText = $nextPara.Range.Text
            $nextTrimmed = $nextText.Trim([char]13, [char]7)

            if ($nextTrimmed -eq "") {
                $nextPara.Style = $joinStyleName
            }
        }
    }
}
